# First blog post for the "news" sheet.
#
# The sheet already had a single placeholder row (id=1, title="Test",
# excerpt="This is a test", category="Test", date="May 08, 2025",
# md_file_name="hello_world.md", img_file_name="hello_word.png").
# Replace the placeholder title/excerpt/category with the real first post,
# and fix the typo'd image file name. The date and the markdown file name
# are already correct, so they are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("news")

# Fix the typo in the image file name first (hello_word.png -> hello_world.png).
$ws.Range("G2").Value = "hello_world.png"

# Fill in the real excerpt, title and category for the announcement post.
$ws.Range("C2").Value = "We’re thrilled to join Purdue University and to unveil our redesigned online home."
$ws.Range("B2").Value = "A New Chapter: Ciampitti Lab Moves to Purdue & Launches a Fresh Website"
$ws.Range("E2").Value = "General"

# Leave the selection where the author ended up after editing the row.
$ws.Range("E3").Select()
